{"js": "// The original document starts with two paragraphs:\n//   1) a \"Title\" styled paragraph containing the text \"Tutorial\"\n//      (center aligned)\n//   2) an otherwise empty paragraph that only carries the\n//      \"_GoBack\" bookmark\n// This edit removes the \"Tutorial\" title paragraph entirely and moves\n// its center alignment onto the paragraph that carried the bookmark,\n// so the document now begins with a single, empty, centered paragraph\n// that still owns the \"_GoBack\" bookmark.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\n// Locate the \"Title\" styled paragraph (currently holding \"Tutorial\").\nlet titleParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.style === \"Title\") {\n    titleParagraph = p;\n    break;\n  }\n}\n\n// The paragraph right after it is the one carrying the \"_GoBack\"\n// bookmark; it inherits the centered alignment.\nconst bookmarkParagraph = titleParagraph.getNext();\nbookmarkParagraph.alignment = Word.Alignment.centered;\n\n// Remove the old \"Title\" paragraph (style, run and paragraph mark).\ntitleParagraph.delete();\n\nawait context.sync();\n", "ps1": "# The document starts with two paragraphs:\n#   1) a \"Title\" styled paragraph containing the text \"Tutorial\"\n#      (center aligned)\n#   2) an otherwise empty paragraph that only carries the\n#      \"_GoBack\" bookmark\n# This edit removes the \"Tutorial\" title paragraph entirely and moves\n# its center alignment onto the paragraph that carried the bookmark,\n# so the document now begins with a single, empty, centered paragraph\n# that still owns the \"_GoBack\" bookmark.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Title\" styled paragraph (currently holding \"Tutorial\").\n$titlePara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Title\") {\n        $titlePara = $p\n        break\n    }\n}\n\n# The paragraph right after it is the one carrying the \"_GoBack\"\n# bookmark; it inherits the centered alignment.\n$bookmarkPara = $titlePara.Next()\n$bookmarkPara.Alignment = 1  # wdAlignParagraphCenter\n\n# Remove the old \"Title\" paragraph (style, run and paragraph mark).\n$titlePara.Range.Delete()\n"}
